{"js": "// Apply left-alignment + single underline to the three Jinja-template\n// paragraphs that make up the `{% for paso in paso_prueba %} ... {% endfor %}`\n// block (\"{% for paso in paso_prueba %}\", \"{{ loop.index }}. {{ paso }}\",\n// \"{% endfor %}\"), matching the authored diff.\n\nconst body = context.document.body;\n\n// Each search string is unique within the document and falls inside one of\n// the three target paragraphs, so expanding to the enclosing paragraph gives\n// us exactly the paragraph we need to restyle.\nconst anchors = [\"paso_prueba\", \"loop.index\", \"endfor\"];\n\nconst paragraphs = [];\nfor (const needle of anchors) {\n  const results = body.search(needle, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find anchor text: \" + needle);\n  }\n\n  const para = results.items[0].paragraphs.getFirst();\n  paragraphs.push(para);\n}\n\nfor (const para of paragraphs) {\n  para.alignment = Word.Alignment.left;\n  para.font.underline = Word.UnderlineType.single;\n}\n\nawait context.sync();\n", "ps1": "# Apply left-alignment + single underline to the three Jinja-template\n# paragraphs that make up the `{% for paso in paso_prueba %} ... {% endfor %}`\n# block (\"{% for paso in paso_prueba %}\", \"{{ loop.index }}. {{ paso }}\",\n# \"{% endfor %}\"), matching the authored diff.\n\n$d = $word.ActiveDocument\n\n$wdLeft = [Microsoft.Office.Interop.Word.WdParagraphAlignment]::wdAlignParagraphLeft\n$wdSingleUnderline = [Microsoft.Office.Interop.Word.WdUnderline]::wdUnderlineSingle\n\n$anchors = \"paso_prueba\", \"loop.index\", \"endfor\"\n\nforeach ($needle in $anchors) {\n    $rng = $d.Content\n    $rng.Find.Text = $needle\n    $rng.Find.Execute() | Out-Null\n\n    $para = $rng.Paragraphs(1)\n    $pRange = $para.Range\n\n    $pRange.ParagraphFormat.Alignment = $wdLeft\n    $pRange.Font.Underline = $wdSingleUnderline\n}\n"}
